# [ADDITIONAL SCRAPING] added code to scrape more data about a player's
# batting performance in a match, also updated the excel sheets.
#
# - Adds a new "Player Info" sheet at the front.
# - Renames MATCH_CARD_LINK -> MATCH_CODE on "ODI Batting" / "ODI Bowling"
#   and rewrites the cell values from the full howstat URL down to just the
#   numeric match code.
# - Drops the now-pointless empty INNING_NUMBER cells on "ODI Batting" rows
#   where the player did not bat.
# - Adds a new "ODI Batting Extra" sheet at the end with additional
#   per-innings batting detail for the more recent matches.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) New "Player Info" sheet, inserted before the existing first sheet.
# ---------------------------------------------------------------------
$first = $wb.Worksheets.Item(1)
$playerInfo = $wb.Worksheets.Add($first)
$playerInfo.Name = "Player Info"

$piHeaders = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($c = 1; $c -le $piHeaders.Length; $c++) {
    $playerInfo.Cells.Item(1, $c).Value = $piHeaders[$c - 1]
}
$piHeaderRange = $playerInfo.Range("A1:D1")
$piHeaderRange.Font.Bold = $true
$piHeaderRange.Borders.LineStyle = 1
$piHeaderRange.HorizontalAlignment = -4108
$piHeaderRange.VerticalAlignment = -4160

$playerInfo.Cells.Item(2, 1).Value = "'3936"
$playerInfo.Cells.Item(2, 2).Value = "Colin de Grandhomme"
$playerInfo.Cells.Item(2, 3).Value = "Right Handed"
$playerInfo.Cells.Item(2, 4).Value = "Right Arm Fast Medium"

# ---------------------------------------------------------------------
# 2) "ODI Batting": MATCH_CARD_LINK -> MATCH_CODE, URL -> bare code,
#    drop empty INNING_NUMBER cells for matches not batted in.
# ---------------------------------------------------------------------
$batting = $wb.Worksheets.Item("ODI Batting")
$batting.Cells.Item(1, 4).Value = "MATCH_CODE"

$lastRow = $batting.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $batting.Cells.Item($r, 4)
    $url = $cell.Value2
    $code = $url -replace ".*MatchCode=", ""
    $cell.Value = "'" + $code
}

$emptyInningRows = @(7, 9, 11, 23, 25, 27, 28, 30, 32, 44)
foreach ($r in $emptyInningRows) {
    $batting.Cells.Item($r, 2).ClearContents()
}

# ---------------------------------------------------------------------
# 3) "ODI Bowling": MATCH_CARD_LINK -> MATCH_CODE, URL -> bare code.
# ---------------------------------------------------------------------
$bowling = $wb.Worksheets.Item("ODI Bowling")
$bowling.Cells.Item(1, 2).Value = "MATCH_CODE"

$lastRow = $bowling.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $bowling.Cells.Item($r, 2)
    $url = $cell.Value2
    $code = $url -replace ".*MatchCode=", ""
    $cell.Value = "'" + $code
}

# ---------------------------------------------------------------------
# 4) New "ODI Batting Extra" sheet, appended after "ODI Bowling".
# ---------------------------------------------------------------------
$bowling = $wb.Worksheets.Item("ODI Bowling")
$extra = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $bowling)
$extra.Name = "ODI Batting Extra"

$exHeaders = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($c = 1; $c -le $exHeaders.Length; $c++) {
    $extra.Cells.Item(1, $c).Value = $exHeaders[$c - 1]
}
$exHeaderRange = $extra.Range("A1:F1")
$exHeaderRange.Font.Bold = $true
$exHeaderRange.Borders.LineStyle = 1
$exHeaderRange.HorizontalAlignment = -4108
$exHeaderRange.VerticalAlignment = -4160

# MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
$extraRows = @(
    @("4250", 6, $null, $null, $null, "NO"),
    @("4251", 6, $null, $null, $null, "NO"),
    @("4252", 7, "4", "2", "11.21%", "NO"),
    @("4305", 7, $null, $null, $null, "NO"),
    @("4311", $null, $null, $null, $null, "NO"),
    @("4315", 7, $null, $null, $null, "NO"),
    @("4328", 7, "5", "2", "24.49%", "NO"),
    @("4333", $null, $null, $null, $null, "NO"),
    @("4337", $null, $null, $null, $null, "NO"),
    @("4341", 6, "0", "0", $null, "NO"),
    @("4346", $null, $null, $null, $null, "NO"),
    @("4353", $null, $null, $null, $null, "NO"),
    @("4355", $null, $null, $null, $null, "NO"),
    @("4402", $null, $null, $null, $null, "NO"),
    @("4406", $null, $null, $null, $null, "NO"),
    @("4410", 7, "6", "3", "19.33%", "NO"),
    @("4423", 7, "3", "0", "13.37%", "NO"),
    @("4563", 7, $null, $null, $null, "NO"),
    @("4566", 7, "1", "0", "6.06%", "NO"),
    @("4568", $null, $null, $null, $null, "NO")
)

$r = 2
foreach ($row in $extraRows) {
    $extra.Cells.Item($r, 1).Value = "'" + $row[0]
    if ($null -ne $row[1]) {
        $extra.Cells.Item($r, 2).Value = $row[1]
    }
    if ($null -ne $row[2]) {
        $extra.Cells.Item($r, 3).Value = "'" + $row[2]
    }
    if ($null -ne $row[3]) {
        $extra.Cells.Item($r, 4).Value = "'" + $row[3]
    }
    if ($null -ne $row[4]) {
        $extra.Cells.Item($r, 5).Value = "'" + $row[4]
    }
    $extra.Cells.Item($r, 6).Value = $row[5]
    $r++
}
